# Update the dSF column (F) values for the bumgarner_madison data sheet.
# These values represent "repulled" data per the commit message
# ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = -4
    4  = -3
    5  = -1
    6  = -7
    8  = -5
    9  = -3
    10 = -3
    11 = -3
    12 = 5
    13 = -2
    14 = 2
    15 = 1
    16 = -1
    17 = -5
    18 = 5
    19 = -3
    20 = -4
    21 = -2
    22 = -4
    23 = -3
    24 = -1
    25 = 4
    26 = -1
    27 = 1
    29 = -1
    30 = -1
    31 = -1
    32 = 5
    33 = 1
    34 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
